# format socket stream as per configuration
#
# Adds a "datatype" column (K) describing each parameter row's data type,
# rewrites the timestamp row's reference-Julian-date example values, and
# gives those example cells (G2:H2) a two-decimal numeric format instead
# of scientific notation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "datatype" column header (matches the bold style already used
#     by the other plain header cells such as A1/I1/J1) ---
$ws.Range("K1").Value = "datatype"
$ws.Range("K1").Font.Bold = $true

# --- Update the sample Julian-date numbers for the "timestamp" row and
#     re-label the notes cell accordingly ---
$ws.Range("D2").Value = 2415020.5
$ws.Range("E2").Value = 2444239.5
$ws.Range("G2").Value = 2469807.5
$ws.Range("H2").Value = 2524593.5
$ws.Range("G2:H2").NumberFormat = "0.00;[Red]0.00"
$ws.Range("J2").Value = "Julian date"

# --- Populate the new "datatype" column for every parameter row ---
$ws.Range("K2").Value = "date"
$ws.Range("K3").Value = "number"
$ws.Range("K4").Value = "number"
$ws.Range("K5").Value = "number"
$ws.Range("K6").Value = "number"
$ws.Range("K7").Value = "number"
$ws.Range("K8").Value = "number"
$ws.Range("K9").Value = "number"
$ws.Range("K10").Value = "number"
$ws.Range("K11").Value = "string"

# --- Widen column G slightly to fit the new values, and move the
#     active selection to D4 ---
$ws.Columns.Item(7).ColumnWidth = 10.998697916666666
$ws.Range("D4").Select() | Out-Null
